{"js": "// Update the date line and every \"a\u00f7b=\" problem cell in the practice\n// table (each row's worksheet was regenerated, so both the operands and\n// the quotient-divisor pair change). Replacements are applied in the\n// same order the text appears in the document; when the same \"before\"\n// text occurs more than once (e.g. \"678\u00f76=\" shows up twice), each call\n// below resolves to that occurrence's own distinct \"after\" text because\n// every search re-scans the (already partially edited) body and always\n// targets the first remaining match.\nconst replacements = [\n  [\"2024-12-13 Friday\", \"2024-12-14 Saturday\"],\n  [\"326\u00f74=\", \"937\u00f79=\"],\n  [\"680\u00f74=\", \"272\u00f75=\"],\n  [\"189\u00f77=\", \"598\u00f78=\"],\n  [\"490\u00f78=\", \"432\u00f76=\"],\n  [\"872\u00f72=\", \"342\u00f79=\"],\n  [\"690\u00f78=\", \"707\u00f74=\"],\n  [\"759\u00f74=\", \"253\u00f79=\"],\n  [\"159\u00f72=\", \"681\u00f77=\"],\n  [\"678\u00f76=\", \"845\u00f73=\"], // first \"678\u00f76=\" occurrence\n  [\"118\u00f72=\", \"116\u00f76=\"],\n  [\"864\u00f75=\", \"378\u00f78=\"],\n  [\"819\u00f79=\", \"427\u00f73=\"],\n  [\"667\u00f77=\", \"646\u00f78=\"],\n  [\"498\u00f76=\", \"332\u00f79=\"],\n  [\"816\u00f74=\", \"448\u00f74=\"],\n  [\"302\u00f75=\", \"129\u00f77=\"],\n  [\"185\u00f77=\", \"374\u00f72=\"],\n  [\"891\u00f74=\", \"148\u00f73=\"],\n  [\"974\u00f78=\", \"100\u00f72=\"],\n  [\"493\u00f76=\", \"940\u00f75=\"],\n  [\"678\u00f76=\", \"241\u00f75=\"], // second \"678\u00f76=\" occurrence\n  [\"485\u00f73=\", \"951\u00f75=\"],\n  [\"610\u00f72=\", \"454\u00f74=\"],\n  [\"823\u00f73=\", \"194\u00f77=\"],\n  [\"115\u00f75=\", \"239\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n  // Replace only the first still-unreplaced match, preserving its\n  // existing run formatting (font/size) since insertText with\n  // \"Replace\" swaps the text inside the matched range in place.\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"a\u00f7b=\" problem cell in the practice\n# table (each row's worksheet was regenerated, so both the operands and\n# the quotient-divisor pair change). Pairs are listed in the same order\n# the text appears in the document; when the same \"before\" text occurs\n# more than once (e.g. \"678\u00f76=\" shows up twice), each loop iteration\n# still resolves to that occurrence's own distinct \"after\" text because\n# Find.Execute with wdReplaceOne re-scans $d.Content from the top and\n# always stops at the first remaining (i.e. not yet replaced) match.\n$d = $word.ActiveDocument\n$wdReplaceOne = 1\n\n$pairs = @(\n    ,@(\"2024-12-13 Friday\", \"2024-12-14 Saturday\")\n    ,@(\"326\u00f74=\", \"937\u00f79=\")\n    ,@(\"680\u00f74=\", \"272\u00f75=\")\n    ,@(\"189\u00f77=\", \"598\u00f78=\")\n    ,@(\"490\u00f78=\", \"432\u00f76=\")\n    ,@(\"872\u00f72=\", \"342\u00f79=\")\n    ,@(\"690\u00f78=\", \"707\u00f74=\")\n    ,@(\"759\u00f74=\", \"253\u00f79=\")\n    ,@(\"159\u00f72=\", \"681\u00f77=\")\n    ,@(\"678\u00f76=\", \"845\u00f73=\")   # first \"678\u00f76=\" occurrence\n    ,@(\"118\u00f72=\", \"116\u00f76=\")\n    ,@(\"864\u00f75=\", \"378\u00f78=\")\n    ,@(\"819\u00f79=\", \"427\u00f73=\")\n    ,@(\"667\u00f77=\", \"646\u00f78=\")\n    ,@(\"498\u00f76=\", \"332\u00f79=\")\n    ,@(\"816\u00f74=\", \"448\u00f74=\")\n    ,@(\"302\u00f75=\", \"129\u00f77=\")\n    ,@(\"185\u00f77=\", \"374\u00f72=\")\n    ,@(\"891\u00f74=\", \"148\u00f73=\")\n    ,@(\"974\u00f78=\", \"100\u00f72=\")\n    ,@(\"493\u00f76=\", \"940\u00f75=\")\n    ,@(\"678\u00f76=\", \"241\u00f75=\")   # second \"678\u00f76=\" occurrence\n    ,@(\"485\u00f73=\", \"951\u00f75=\")\n    ,@(\"610\u00f72=\", \"454\u00f74=\")\n    ,@(\"823\u00f73=\", \"194\u00f77=\")\n    ,@(\"115\u00f75=\", \"239\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceOne) | Out-Null\n}\n"}
